$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("0.06344039738178253", "0.9872387051582336", "0.01868856512010098", "0.9983006715774536")
    ,@("0.009622324258089066", "0.9982057809829712", "0.009350425563752651", "0.9986274242401123")
    ,@("0.004750330466777086", "0.9986991882324219", "0.04307442530989647", "0.9842483401298523")
    ,@("0.002519248286262155", "0.9993944764137268", "0.1015850231051445", "0.9587581753730774")
    ,@("0.001754428492859006", "0.9993944764137268", "0.005921387579292059", "0.998431384563446")
    ,@("0.001033064094372094", "0.9996860027313232", "0.002392819384112954", "0.9995424747467041")
    ,@("0.000368318083928898", "0.9999327063560486", "0.00712552247568965", "0.9970588088035583")
    ,@("0.001066728611476719", "0.9995738863945007", "0.003266300074756145", "0.9990849494934082")
    ,@("0.001054296619258821", "0.9996635913848877", "0.003284576814621687", "0.9996078610420227")
    ,@("0.0003046854108106345", "0.9999327063560486", "0.001471647177822888", "0.999803900718689")
    ,@("0.0007782182656228542", "0.9998205900192261", "0.002296320861205459", "0.9997385740280151")
    ,@("0.001022602431476116", "0.9996635913848877", "0.0016708456678316", "0.9996731877326965")
    ,@("8.251862163888291E-05", "1", "0.00180995324626565", "0.9996078610420227")
    ,@("0.0007882217178121209", "0.9998430013656616", "0.001503518666140735", "0.9997385740280151")
    ,@("0.0005129491328261793", "0.9997981786727905", "0.0008394100586883724", "0.9998692870140076")
    ,@("0.0003737210354302078", "0.9998430013656616", "0.02692805044353008", "0.9896078705787659")
    ,@("0.000610753137152642", "0.9998654127120972", "0.001625268836505711", "0.999803900718689")
    ,@("0.0003747462760657072", "0.9998430013656616", "0.001675666542723775", "0.9998692870140076")
    ,@("0.0001375305437250063", "0.9999551177024841", "0.003317722352221608", "0.999803900718689")
    ,@("0.0003153570287395269", "0.9999327063560486", "0.002282301662489772", "0.9998692870140076")
    ,@("0.0008108025649562478", "0.9998654127120972", "0.002164163626730442", "0.9998692870140076")
    ,@("0.0002270948025397956", "0.999910295009613", "0.006921723484992981", "0.9997385740280151")
    ,@("0.0003711171739269048", "0.9999551177024841", "0.002539742272347212", "0.9998692870140076")
    ,@("0.001101163681596518", "0.9997981786727905", "0.003056906629353762", "0.9998692870140076")
    ,@("0.0002077439421555027", "0.999910295009613", "0.001975724240764976", "0.9996731877326965")
    ,@("0.0001041342038661242", "0.9999551177024841", "0.002347868634387851", "0.9998692870140076")
    ,@("0.0002792274754028767", "0.999910295009613", "0.002295139012858272", "0.9998692870140076")
    ,@("0.0004129930457565933", "0.9998654127120972", "0.003743712091818452", "0.9998692870140076")
    ,@("0.0003320025571156293", "0.9998878836631775", "0.003356365719810128", "0.9998692870140076")
    ,@("6.80771263432689E-05", "0.9999775886535645", "0.003379174042493105", "0.9998692870140076")
    ,@("0.0002117311232723296", "0.9999775886535645", "0.004439869895577431", "0.9995424747467041")
    ,@("0.000303563429042697", "0.9999775886535645", "0.008532056584954262", "0.9981045722961426")
    ,@("0.0001228184555657208", "0.9999775886535645", "0.004903410095721483", "0.9996078610420227")
    ,@("0.001190558425150812", "0.9998430013656616", "0.002360643586143851", "0.999803900718689")
    ,@("0.0001502480590716004", "0.9999775886535645", "0.009292657487094402", "0.9958823323249817")
    ,@("0.0005973048973828554", "0.9998430013656616", "0.003298654919490218", "0.9998692870140076")
    ,@("0.0001404710201313719", "0.9999775886535645", "0.004033519886434078", "0.9997385740280151")
    ,@("0.0001499985373811796", "0.9999551177024841", "0.2056060433387756", "0.973006546497345")
    ,@("0.0005677217268384993", "0.9998878836631775", "0.002681964309886098", "0.9996731877326965")
    ,@("4.17267219745554E-05", "0.9999775886535645", "0.002199480077251792", "0.999803900718689")
    ,@("0.0002282280765939504", "0.9999775886535645", "0.003055559238418937", "0.9996731877326965")
    ,@("9.227167902281508E-05", "0.9999775886535645", "0.004759395495057106", "0.9996078610420227")
    ,@("0.0008515903609804809", "0.9998654127120972", "0.005835357587784529", "0.9994117617607117")
    ,@("0.0002261419431306422", "0.9999775886535645", "0.004704890307039022", "0.9994771480560303")
    ,@("0.0001170756586361676", "0.9999775886535645", "0.003880236996337771", "0.9994771480560303")
    ,@("2.512126229703426E-05", "1", "0.004128599539399147", "0.9996078610420227")
    ,@("0.0003317774971947074", "0.9999551177024841", "0.003538544289767742", "0.9996078610420227")
    ,@("0.0004645197768695652", "0.999910295009613", "0.008834082633256912", "0.9969934821128845")
    ,@("0.0001582972618052736", "0.9999775886535645", "0.003768350463360548", "0.9998692870140076")
    ,@("0.0007153008482418954", "0.9999551177024841", "0.01460480503737926", "0.9946405291557312")
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value2 = [double]$row[0]
    $ws.Cells.Item($r, 2).Value2 = [double]$row[1]
    $ws.Cells.Item($r, 3).Value2 = [double]$row[2]
    $ws.Cells.Item($r, 4).Value2 = [double]$row[3]
}
